$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 109.875
$ws.Range("I6").Value = 68.42856999999999
$ws.Range("K6").Value = 205.28571
$ws.Range("M6").Value = -93.28570999999999
$ws.Range("H37").Value = 98
$ws.Range("I37").Value = 98
$ws.Range("K37").Value = 294
$ws.Range("M37").Value = -168
$ws.Range("H42").Value = 114.75
$ws.Range("I42").Value = 86.666664
$ws.Range("J42").Value = 199
$ws.Range("K42").Value = 259.999992
$ws.Range("L42").Value = 597
$ws.Range("M42").Value = -29.99999200000002
$ws.Range("N42").Value = -1057
$ws.Range("H94").Value = 802
$ws.Range("I94").Value = 504
$ws.Range("K94").Value = 504
$ws.Range("M94").Value = -53

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 1375
$ws.Range("I21").Value = 750
$ws.Range("J21").Value = 2000
$ws.Range("K21").Value = 750
$ws.Range("L21").Value = 2000
$ws.Range("M21").Value = -376
$ws.Range("N21").Value = -2748
$ws.Range("H32").Value = 1119.5834
$ws.Range("I32").Value = 1117.2858
$ws.Range("K32").Value = 1117.2858
$ws.Range("M32").Value = -830.2858000000001
$ws.Range("H37").Value = 47008.25
$ws.Range("J37").Value = 61666.332
$ws.Range("L37").Value = 61666.332
$ws.Range("N37").Value = -62212.332
$ws.Range("H45").Value = 1775.625
$ws.Range("I45").Value = 1775.625
$ws.Range("K45").Value = 1775.625
$ws.Range("M45").Value = -1398.625
$ws.Range("H110").Value = 10057.9
$ws.Range("I110").Value = 10620
$ws.Range("K110").Value = 10620
$ws.Range("M110").Value = -8575
$ws.Range("H134").Value = 10000
$ws.Range("J134").Value = 10000
$ws.Range("L134").Value = 10000
$ws.Range("N134").Value = -20140
$ws.Range("H141").Value = 15000
$ws.Range("J141").Value = 15000
$ws.Range("L141").Value = 15000
$ws.Range("N141").Value = -25360

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1050
$ws.Range("I107").Value = 900
$ws.Range("K107").Value = 900
$ws.Range("M107").Value = 1020
$ws.Range("H134").Value = 1791.0834
$ws.Range("I134").Value = 1537.9565
$ws.Range("K134").Value = 4613.8695
$ws.Range("M134").Value = -2078.8695

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1332.44
$ws.Range("I31").Value = 1035.8
$ws.Range("J31").Value = 1777.4
$ws.Range("K31").Value = 1035.8
$ws.Range("L31").Value = 1777.4
$ws.Range("M31").Value = -740.8
$ws.Range("N31").Value = -2367.4
$ws.Range("H34").Value = 1332.44
$ws.Range("I34").Value = 1035.8
$ws.Range("J34").Value = 1777.4
$ws.Range("K34").Value = 1035.8
$ws.Range("L34").Value = 1777.4
$ws.Range("M34").Value = -833.8
$ws.Range("N34").Value = -2181.4
$ws.Range("H58").Value = 2250
$ws.Range("I58").Value = 2250
$ws.Range("K58").Value = 2250
$ws.Range("M58").Value = -2047
$ws.Range("H62").Value = 7630.3335
$ws.Range("I62").Value = 7599
$ws.Range("J62").Value = 7661.6665
$ws.Range("K62").Value = 7599
$ws.Range("L62").Value = 7661.6665
$ws.Range("M62").Value = -6975
$ws.Range("N62").Value = -8909.666499999999
$ws.Range("H65").Value = 7630.3335
$ws.Range("I65").Value = 7599
$ws.Range("J65").Value = 7661.6665
$ws.Range("K65").Value = 37995
$ws.Range("L65").Value = 38308.3325
$ws.Range("M65").Value = -34875
$ws.Range("N65").Value = -44548.3325
$ws.Range("H122").Value = 2799.0667
$ws.Range("I122").Value = 2799.0667
$ws.Range("K122").Value = 8397.2001
$ws.Range("M122").Value = -5947.2001
$ws.Range("H132").Value = 2489.75
$ws.Range("I132").Value = 2738
$ws.Range("K132").Value = 8214
$ws.Range("M132").Value = -5684
$ws.Range("H134").Value = 2195.8147
$ws.Range("I134").Value = 1931.52
$ws.Range("K134").Value = 5794.559999999999
$ws.Range("M134").Value = -3259.559999999999
$ws.Range("H136").Value = 2250
$ws.Range("I136").Value = 2250
$ws.Range("K136").Value = 6750
$ws.Range("M136").Value = -4200

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1152283.9
$ws.Range("I4").Value = 674935.8
$ws.Range("J4").Value = 4289143
$ws.Range("K4").Value = 2024807.4
$ws.Range("L4").Value = 12867429
$ws.Range("M4").Value = -2024695.4
$ws.Range("N4").Value = -12867653
$ws.Range("H8").Value = 3011.7144
$ws.Range("I8").Value = 3011.7144
$ws.Range("K8").Value = 9035.143199999999
$ws.Range("M8").Value = -8896.143199999999
$ws.Range("H11").Value = 1132.5
$ws.Range("I11").Value = 1132.5
$ws.Range("K11").Value = 3397.5
$ws.Range("M11").Value = -3257.5
$ws.Range("H68").Value = 2691.0715
$ws.Range("J68").Value = 2928.15
$ws.Range("L68").Value = 8784.450000000001
$ws.Range("N68").Value = -10406.45
$ws.Range("H70").Value = 13349
$ws.Range("I70").Value = 13349
$ws.Range("K70").Value = 40047
$ws.Range("M70").Value = -39732
$ws.Range("H71").Value = 2691.0715
$ws.Range("J71").Value = 2928.15
$ws.Range("L71").Value = 26353.35
$ws.Range("N71").Value = -34465.35000000001
$ws.Range("H73").Value = 13349
$ws.Range("I73").Value = 13349
$ws.Range("K73").Value = 40047
$ws.Range("M73").Value = -38955
$ws.Range("H75").Value = 2271.4285
$ws.Range("I75").Value = 2129.6667
$ws.Range("K75").Value = 6389.000100000001
$ws.Range("M75").Value = -5391.000100000001
$ws.Range("H78").Value = 2271.4285
$ws.Range("I78").Value = 2129.6667
$ws.Range("K78").Value = 19167.0003
$ws.Range("M78").Value = -14175.0003
$ws.Range("H113").Value = 537.8570999999999
$ws.Range("J113").Value = 619.6
$ws.Range("L113").Value = 1858.8
$ws.Range("N113").Value = -6198.8
$ws.Range("H137").Value = 1030
$ws.Range("I137").Value = 1030
$ws.Range("K137").Value = 3090
$ws.Range("M137").Value = 2010

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 101.90909
$ws.Range("I2").Value = 98.3125
$ws.Range("J2").Value = 111.5
$ws.Range("K2").Value = 98.3125
$ws.Range("L2").Value = 111.5
$ws.Range("M2").Value = 14.6875
$ws.Range("N2").Value = -337.5
$ws.Range("H43").Value = 22112.8
$ws.Range("J43").Value = 35555
$ws.Range("L43").Value = 35555
$ws.Range("N43").Value = -35857
$ws.Range("H139").Value = 25000
$ws.Range("J139").Value = 25000
$ws.Range("L139").Value = 25000
$ws.Range("N139").Value = -35280

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 2499
$ws.Range("I22").Value = 2499
$ws.Range("K22").Value = 2499
$ws.Range("M22").Value = -2204
$ws.Range("H27").Value = 2499
$ws.Range("I27").Value = 2499
$ws.Range("K27").Value = 2499
$ws.Range("M27").Value = -2392
$ws.Range("H46").Value = 2688.125
$ws.Range("I46").Value = 2294.1428
$ws.Range("J46").Value = 2994.5557
$ws.Range("K46").Value = 2294.1428
$ws.Range("L46").Value = 2994.5557
$ws.Range("M46").Value = -2106.1428
$ws.Range("N46").Value = -3370.5557
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
$ws.Range("H132").Value = 4055.375
$ws.Range("I132").Value = 2740
$ws.Range("J132").Value = 8001.5
$ws.Range("K132").Value = 8220
$ws.Range("L132").Value = 24004.5
$ws.Range("M132").Value = -5690
$ws.Range("N132").Value = -29064.5
$ws.Range("H136").Value = 55557500
$ws.Range("I136").Value = 1758.6
$ws.Range("K136").Value = 5275.799999999999
$ws.Range("M136").Value = -2725.799999999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H122").Value = 2914.111
$ws.Range("I122").Value = 2914.111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8742.332999999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6292.332999999999
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 1995.6
$ws.Range("I126").Value = 1995.6
$ws.Range("K126").Value = 5986.799999999999
$ws.Range("M126").Value = -3516.799999999999
$ws.Range("H132").Value = 7354.2856
$ws.Range("I132").Value = 5695
$ws.Range("K132").Value = 17085
$ws.Range("M132").Value = -14555
